# Refresh the cryptos price/volume snapshot (Price = column D, Volume(1h) = column E).
# Price cells that look like plain numbers ("213.76", "33.00", ...) are forced to
# text via NumberFormat "@" before assignment so Excel doesn't silently coerce them
# to doubles (which would lose trailing zeros / change "."-grouped values like
# "3.226.48"); the style is reset back to Normal afterwards so no visible formatting
# change is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '80.805.51'
$ws.Range("E2").Value = '  +5.86%  '
$ws.Range("D3").Value = '3.228.70'
$ws.Range("E3").Value = '  +6.03%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.76'
$ws.Range("E5").Value = '  +7.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.14'
$ws.Range("E6").Value = '  +2.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.281'
$ws.Range("E7").Value = '  +34.86%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.609'
$ws.Range("E9").Value = '  +10.60%  '
$ws.Range("D10").Value = '3.226.48'
$ws.Range("E10").Value = '  +6.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.620'
$ws.Range("E11").Value = '  +41.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000277'
$ws.Range("E12").Value = '  +43.81%  '
$ws.Range("E13").Value = '  +3.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.48'
$ws.Range("E14").Value = '  +5.38%  '
$ws.Range("D15").Value = '3.820.97'
$ws.Range("E15").Value = '  +6.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '33.00'
$ws.Range("E16").Value = '  +14.25%  '
$ws.Range("D17").Value = '80.743.73'
$ws.Range("E17").Value = '  +5.95%  '
$ws.Range("D18").Value = '3.217.50'
$ws.Range("E18").Value = '  +5.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.68'
$ws.Range("E19").Value = '  +8.52%  '
$ws.Range("E20").Value = '  +26.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.42'
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '449.69'
$ws.Range("E22").Value = '  +17.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.43'
$ws.Range("E23").Value = '  +23.71%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.88'
$ws.Range("E24").Value = '  +12.62%  '
$ws.Range("D25").Value = '3.389.16'
$ws.Range("E25").Value = '  +6.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '78.02'
$ws.Range("E26").Value = '  +7.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.12'
$ws.Range("E27").Value = '  +13.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000129'
$ws.Range("E28").Value = '  +19.34%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.43'
$ws.Range("E30").Value = '  +13.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '569.78'
$ws.Range("E32").Value = '  +15.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.53'
$ws.Range("E33").Value = '  +9.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.161'
$ws.Range("E34").Value = '  +31.42%  '
$ws.Range("E35").Value = '  +7.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.86'
$ws.Range("E36").Value = '  +15.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.125'
$ws.Range("E37").Value = '  +21.13%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.420'
$ws.Range("E38").Value = '  +11.25%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.86'
$ws.Range("E40").Value = '  +14.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '164.28'
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.37'
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '193.27'
$ws.Range("E43").Value = '  +0.56%  '
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("E45").Value = '  +13.08%  '
$ws.Range("E46").Value = '  +13.67%  '
$ws.Range("E47").Value = '  +9.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.805'
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("E49").Value = '  +13.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.66'
$ws.Range("E50").Value = '  +6.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.651'
$ws.Range("E51").Value = '  +9.64%  '

# Restore default (unstyled) formatting on price cells forced to text above
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
